$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 270 (existing rows 270..295 shift down to 272..297)
$ws.Rows("270:271").Insert()

# Row 270 - new record (Coliflor, Primera)
$ws.Range("A270").Value = 5
$ws.Range("B270").Value = "Macroferia Regional de Talca"
$ws.Range("C270").Value = "Maule"
$ws.Range("D270").Value = 44783
$ws.Range("E270").Value = 7
$ws.Range("F270").Value = 100112008
$ws.Range("G270").Value = "Coliflor"
$ws.Range("H270").Value = "Sin especificar"
$ws.Range("I270").Value = "Primera"
$ws.Range("J270").Value = 3000
$ws.Range("K270").Value = 1000
$ws.Range("L270").Value = 1000
$ws.Range("M270").Value = 1000
$ws.Range("N270").Value = "`$/unidad"
$ws.Range("O270").Value = "Región del Maule"
$ws.Range("P270").Value = 1000
$ws.Range("Q270").Value = 1
$ws.Range("R270").Value = "Hortaliza"

# Row 271 - new record (Coliflor, Segunda)
$ws.Range("A271").Value = 5
$ws.Range("B271").Value = "Macroferia Regional de Talca"
$ws.Range("C271").Value = "Maule"
$ws.Range("D271").Value = 44783
$ws.Range("E271").Value = 7
$ws.Range("F271").Value = 100112008
$ws.Range("G271").Value = "Coliflor"
$ws.Range("H271").Value = "Sin especificar"
$ws.Range("I271").Value = "Segunda"
$ws.Range("J271").Value = 2000
$ws.Range("K271").Value = 800
$ws.Range("L271").Value = 800
$ws.Range("M271").Value = 800
$ws.Range("N271").Value = "`$/unidad"
$ws.Range("O271").Value = "Región del Maule"
$ws.Range("P271").Value = 800
$ws.Range("Q271").Value = 1
$ws.Range("R271").Value = "Hortaliza"
